$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate row 8 (date + price text) into a new row 9, copying values & formatting
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4104)
$excel.CutCopyMode = $false
